# limited_offer_instructions: change all mentions of "cents" (the cent sign
# price tags and the explanatory prose) to "points", and insert a new slide
# (position 12) that explains how points convert to cents, pushing the old
# closing "payment" slide from position 12 to position 13.

$p = $ppt.ActivePresentation

# EMU-per-point conversion constant used by the PowerPoint object model
# (Shape.Left/Top/Width/Height are expressed in points).
$EMU = 12700.0

# ---------------------------------------------------------------------
# 1. Update the "12 ¢" price-tag textboxes on slides 3-10 to say
#    "12 points" instead, and nudge each box back to where PowerPoint's
#    autosize (spAutoFit) would have re-centred it once the text grew
#    wider.
# ---------------------------------------------------------------------
$centsBoxes = @(
    @{ Slide = 3;  Shape = 8;  Off = @(5068644, 3302912); Ext = @(2054711, 646331) },
    @{ Slide = 4;  Shape = 8;  Off = @(5068644, 3302912); Ext = @(2054711, 646331) },
    @{ Slide = 5;  Shape = 11; Off = @(5068644, 4014993); Ext = @(2054711, 646331) },
    @{ Slide = 6;  Shape = 7;  Off = @(5068644, 3302912); Ext = @(2054711, 646331) },
    @{ Slide = 7;  Shape = 7;  Off = @(5068644, 3302912); Ext = @(2054711, 646331) },
    @{ Slide = 8;  Shape = 7;  Off = @(5068644, 4163075); Ext = @(2054711, 646331) },
    @{ Slide = 9;  Shape = 8;  Off = @(5068644, 4114282); Ext = @(2054711, 646331) },
    @{ Slide = 10; Shape = 8;  Off = @(5068644, 4114282); Ext = @(2054711, 646331) }
)

foreach ($box in $centsBoxes) {
    $s = $p.Slides.Item($box.Slide)
    $sh = $s.Shapes.Item($box.Shape)
    $sh.TextFrame.TextRange.Text = "12 points"
    $sh.Left = $box.Off[0] / $EMU
    $sh.Top = $box.Off[1] / $EMU
    $sh.Width = $box.Ext[0] / $EMU
    $sh.Height = $box.Ext[1] / $EMU
}

# ---------------------------------------------------------------------
# 2. Update the explanatory sentence on slide 3 that mentions "12 cents"
#    in prose.
# ---------------------------------------------------------------------
$slide3Body = $p.Slides.Item(3).Shapes.Item(5)
$slide3Body.TextFrame.TextRange.Replace("12 cents", "12 points")

# ---------------------------------------------------------------------
# 3. Insert the new "points -> cents" explanation slide at position 12
#    (duplicating slide 11, which carries the right boilerplate Title /
#    Subtitle / Right-Arrow "press the right button" shapes), pushing
#    the old closing slide back to position 13.
# ---------------------------------------------------------------------
$sourceSlide = $p.Slides.Item(11)
$sourceSlide.Duplicate() | Out-Null

$newSlide = $p.Slides.Item(12)
$content = $newSlide.Shapes.Item(4)

$content.Left = 534208 / $EMU
$content.Top = 1218946 / $EMU
$content.Width = 11123583 / $EMU
$content.Height = 3597753 / $EMU

$content.TextFrame.TextRange.Text = "You will receive payment based on the number of points  you win.`rFor every 100 points you win, you will receive 10 cents.`r`r`r100 points = 10 cents"

$para1 = $content.TextFrame.TextRange.Paragraphs(1, 1)
$para1.ParagraphFormat.Alignment = 1

$para5 = $content.TextFrame.TextRange.Paragraphs(5, 1)
$para5.Font.Bold = $true
$para5.Font.Size = 36
